# Regenerate the handback-status report with a new source file pair:
#   355c1993-7e39-430e-b9bb-23a59c13e961.md -> 7adfb1b1-bfa5-40ee-8c09-753191be7360.md
#   c730e75e-4018-49cb-a52d-44a95de63869.md -> ffff9744a076-decc-4a44-a43d-19168c1936af.md
# plus refreshed handoff/handback timestamps and xlf correspondence names.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# ---------------------------------------------------------------------------
# Overview sheet
# ---------------------------------------------------------------------------
$wsOverview.Range("A2").Value = "7adfb1b1-bfa5-40ee-8c09-753191be7360.md"
$wsOverview.Range("B2").Value = "e2e\7adfb1b1-bfa5-40ee-8c09-753191be7360.md"
$wsOverview.Range("G2").Value = "2016-08-23 15:19:57"

$wsOverview.Range("A3").Value = "ffff9744a076-decc-4a44-a43d-19168c1936af.md"
$wsOverview.Range("B3").Value = "e2e\ffff9744a076-decc-4a44-a43d-19168c1936af.md"
$wsOverview.Range("G3").Value = "2016-08-23 15:19:57"

$wsOverview.Hyperlinks.Delete()
$wsOverview.Hyperlinks.Add($wsOverview.Range("B2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/4a5caaae64a0155e5fa3e554a9e6c2b7c4310760/e2e/7adfb1b1-bfa5-40ee-8c09-753191be7360.md", "", "", "e2e\7adfb1b1-bfa5-40ee-8c09-753191be7360.md") | Out-Null
$wsOverview.Hyperlinks.Add($wsOverview.Range("B3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/4a5caaae64a0155e5fa3e554a9e6c2b7c4310760/e2e/ffff9744a076-decc-4a44-a43d-19168c1936af.md", "", "", "e2e\ffff9744a076-decc-4a44-a43d-19168c1936af.md") | Out-Null

# ---------------------------------------------------------------------------
# zh-cn sheet
# ---------------------------------------------------------------------------
$wsZhCn.Range("A2").Value = "7adfb1b1-bfa5-40ee-8c09-753191be7360.md"
$wsZhCn.Range("G2").Value = "7adfb1b1-bfa5-40ee-8c09-753191be7360.7f0181df8325376f5bcb81041e6dd94ad8c69bd0.zh-cn.xlf"
$wsZhCn.Range("H2").Value = "2016-08-23 15:19:52"
$wsZhCn.Range("I2").Value = "7adfb1b1-bfa5-40ee-8c09-753191be7360.md"
$wsZhCn.Range("J2").Value = "7adfb1b1-bfa5-40ee-8c09-753191be7360.7f0181df8325376f5bcb81041e6dd94ad8c69bd0.zh-cn.xlf"
$wsZhCn.Range("K2").Value = "2016-08-23 15:20:44"

$wsZhCn.Range("A3").Value = "ffff9744a076-decc-4a44-a43d-19168c1936af.md"
$wsZhCn.Range("G3").Value = "7adfb1b1-bfa5-40ee-8c09-753191be7360.7f0181df8325376f5bcb81041e6dd94ad8c69bd0.zh-cn.xlf"
$wsZhCn.Range("H3").Value = "2016-08-23 15:19:52"
$wsZhCn.Range("I3").Value = "ffff9744a076-decc-4a44-a43d-19168c1936af.md"
$wsZhCn.Range("J3").Value = "7adfb1b1-bfa5-40ee-8c09-753191be7360.7f0181df8325376f5bcb81041e6dd94ad8c69bd0.zh-cn.xlf"
$wsZhCn.Range("K3").Value = "2016-08-23 15:20:44"

$wsZhCn.Hyperlinks.Delete()
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("A2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/4a5caaae64a0155e5fa3e554a9e6c2b7c4310760/e2e/7adfb1b1-bfa5-40ee-8c09-753191be7360.md", "", "", "7adfb1b1-bfa5-40ee-8c09-753191be7360.md") | Out-Null
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("I2"), "https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/1351b1833c09b8959ab379dce76c561f188c08ac/e2e/7adfb1b1-bfa5-40ee-8c09-753191be7360.md", "", "", "7adfb1b1-bfa5-40ee-8c09-753191be7360.md") | Out-Null
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("A3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/4a5caaae64a0155e5fa3e554a9e6c2b7c4310760/e2e/ffff9744a076-decc-4a44-a43d-19168c1936af.md", "", "", "ffff9744a076-decc-4a44-a43d-19168c1936af.md") | Out-Null
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("I3"), "https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/1351b1833c09b8959ab379dce76c561f188c08ac/e2e/ffff9744a076-decc-4a44-a43d-19168c1936af.md", "", "", "ffff9744a076-decc-4a44-a43d-19168c1936af.md") | Out-Null

# ---------------------------------------------------------------------------
# de-de sheet
# ---------------------------------------------------------------------------
$wsDeDe.Range("A2").Value = "7adfb1b1-bfa5-40ee-8c09-753191be7360.md"
$wsDeDe.Range("G2").Value = "7adfb1b1-bfa5-40ee-8c09-753191be7360.7f0181df8325376f5bcb81041e6dd94ad8c69bd0.de-de.xlf"
$wsDeDe.Range("H2").Value = "2016-08-23 15:19:57"
$wsDeDe.Range("I2").Value = "7adfb1b1-bfa5-40ee-8c09-753191be7360.md"
$wsDeDe.Range("J2").Value = "7adfb1b1-bfa5-40ee-8c09-753191be7360.7f0181df8325376f5bcb81041e6dd94ad8c69bd0.de-de.xlf"
$wsDeDe.Range("K2").Value = "2016-08-23 15:20:53"

$wsDeDe.Range("A3").Value = "ffff9744a076-decc-4a44-a43d-19168c1936af.md"
$wsDeDe.Range("G3").Value = "7adfb1b1-bfa5-40ee-8c09-753191be7360.7f0181df8325376f5bcb81041e6dd94ad8c69bd0.de-de.xlf"
$wsDeDe.Range("H3").Value = "2016-08-23 15:19:57"
$wsDeDe.Range("I3").Value = "ffff9744a076-decc-4a44-a43d-19168c1936af.md"
$wsDeDe.Range("J3").Value = "7adfb1b1-bfa5-40ee-8c09-753191be7360.7f0181df8325376f5bcb81041e6dd94ad8c69bd0.de-de.xlf"
$wsDeDe.Range("K3").Value = "2016-08-23 15:20:53"

$wsDeDe.Hyperlinks.Delete()
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("A2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/4a5caaae64a0155e5fa3e554a9e6c2b7c4310760/e2e/7adfb1b1-bfa5-40ee-8c09-753191be7360.md", "", "", "7adfb1b1-bfa5-40ee-8c09-753191be7360.md") | Out-Null
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("I2"), "https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/12b470e55f9d0e50bdc8f298a795ae3148f733ad/e2e/7adfb1b1-bfa5-40ee-8c09-753191be7360.md", "", "", "7adfb1b1-bfa5-40ee-8c09-753191be7360.md") | Out-Null
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("A3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/4a5caaae64a0155e5fa3e554a9e6c2b7c4310760/e2e/ffff9744a076-decc-4a44-a43d-19168c1936af.md", "", "", "ffff9744a076-decc-4a44-a43d-19168c1936af.md") | Out-Null
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("I3"), "https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/12b470e55f9d0e50bdc8f298a795ae3148f733ad/e2e/ffff9744a076-decc-4a44-a43d-19168c1936af.md", "", "", "ffff9744a076-decc-4a44-a43d-19168c1936af.md") | Out-Null
